$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.74%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.46%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.192"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.93%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07578"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.08%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.343"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.49%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.684"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.30%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9257"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.98%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.82%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1199"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-4.33%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.03%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09024"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.22%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04163"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.37%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.13%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001292"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.96%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005839"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.80%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'3.350"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.14%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3355"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.60%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.620"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'6.72%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1351"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.44%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2811"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.97%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.72%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001270"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.65%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004074"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-2.10%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.29%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02428"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-1.65%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05154"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.37%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007723"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.63%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1299"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.14%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007624"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'11.60%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'51.51%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008215"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.46%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3107"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.86%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006583"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.15%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.29%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.2666"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'29.77%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004201"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'2.47%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.29%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.29%"
$ws.Range("E51").Style = "Normal"
